$wb = $excel.ActiveWorkbook

# --- Driver sheet: update Execution_Flag column (D) values ---
$driver = $wb.Worksheets.Item("Driver")

# D4: "Yes" -> "yes" (keep the existing cell formatting unchanged)
$driver.Range("D4").Value = "yes"
$driver.Range("D2").Copy()
$driver.Range("D4").PasteSpecial(-4122)

# D5, D6: "no" -> "NO" (copy the already-correct D2 cell, value + format,
# so the original cell formatting is preserved)
$driver.Range("D2").Copy()
$driver.Range("D5").PasteSpecial(-4104)
$driver.Range("D6").PasteSpecial(-4104)

$excel.CutCopyMode = $false

# --- Response sheet: update the selected cell ---
$response = $wb.Worksheets.Item("Response")
$response.Select()
$response.Range("I5").Select()

# Re-activate the Driver sheet and set its selection,
# so it remains the active tab as in the source file.
$driver.Select()
$driver.Range("D5").Select()
